# theday09 atributes.xlsx update
# - insert two new columns before column D, pushing the old "my_file pl" /
#   key-name column from D to F
# - add a new "get offer details - key names" column in D with translated
#   (Polish) key labels
# - leave a narrow blank separator column in E
# - re-apply a (slightly) distinct font to the left-hand block (A:D) while
#   the moved-over column (F) keeps the original font
# - fix a stray missing opening quote in the "Faktura" shared string
# - move the active selection to D14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert two new blank columns at D (old D -> F, nothing lands in E)
# ---------------------------------------------------------------------
$ws.Range("D1:E1").EntireColumn.Insert()

# Column widths: D mirrors C's width, E is a narrow spacer, F already
# carries over the old column D width/format because it was shifted.
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth
$ws.Columns("E").ColumnWidth = 1.5

# The Insert() above copies column C's left-aligned look onto every new
# D/E cell in the used range. D2:D18 should end up with no explicit
# alignment (like A/B), and E2:E18 should stay completely untouched/
# empty (only E1 is a real, styled cell) - so wipe both back to a clean
# slate first.
$ws.Range("D2:D18").ClearFormats()
$ws.Range("E2:E18").Clear()

# ---------------------------------------------------------------------
# 2. New header (D1) and translated key-name labels (D2:D18)
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "get offer details - key names"

$ws.Range("D2").Value  = "ID oferty"
$ws.Range("D3").Value  = "ID sprzedajacego"
$ws.Range("D4").Value  = "Lokalizacja"
$ws.Range("D5").Value  = "Tytul"
$ws.Range("D6").Value  = "Cena"
$ws.Range("D7").Value  = "Marka"
$ws.Range("D8").Value  = "Model"
$ws.Range("D9").Value  = "Rok produkcji"
$ws.Range("D10").Value = "Przebieg"
$ws.Range("D11").Value = "Pojemność silnika"
$ws.Range("D12").Value = "Moc"
$ws.Range("D13").Value = "Rodzaj paliwa"
$ws.Range("D14").Value = "Kolor"
$ws.Range("D15").Value = "Uszkodzony"
$ws.Range("D16").Value = "Kraj pochodzenia"
$ws.Range("D17").Value = "Napęd"
$ws.Range("D18").Value = "Liczba miejsc"

# ---------------------------------------------------------------------
# 3. Fix the stray missing opening quote on the "Faktura" entry, now in F2
# ---------------------------------------------------------------------
$ws.Range("F2").Value = '"Faktura": "Wystawiam fakturę VAT",'

# ---------------------------------------------------------------------
# 4. D1 gets the same centered look as B1/C1
# ---------------------------------------------------------------------
$ws.Range("D1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 5. Give the A:D block its own (distinct) font, leaving F with the
#    original font
# ---------------------------------------------------------------------
$ws.Range("A1:D18").Font.ThemeColor = 1

# ---------------------------------------------------------------------
# 6. Update the active selection
# ---------------------------------------------------------------------
$ws.Range("D14").Select()
